$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.931.26'
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("D3").Value = '3.806.01'
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '434.49'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +4.63%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.51'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.65%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.998'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.739'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.25%  '
$ws.Range("E10").Value = '  -11.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0000323'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -16.98%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '42.83'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +3.01%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.52'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +2.97%  '
$ws.Range("D14").Value = '4.403.14'
$ws.Range("E14").Value = '  +0.11%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.00'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("E16").Value = '  -0.39%  '
$ws.Range("D17").Value = '3.833.08'
$ws.Range("E17").Value = '  +0.52%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '19.99'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.15%  '
$ws.Range("E19").Value = '  +3.79%  '
$ws.Range("D20").Value = '66.767.19'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '413.07'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.62%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '14.80'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.24%  '
$ws.Range("E23").Value = '  +4.52%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '85.33'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '37.11'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.53%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.35'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +5.20%  '
$ws.Range("B27").Value = 'LEO'
$ws.Range("C27").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.59'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.38%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.74'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +31.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.88'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.54%  '
$ws.Range("E30").Value = '  +11.22%  '
$ws.Range("B31").Value = 'Bittensor'
$ws.Range("C31").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '714.80'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.10%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.82'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +10.20%  '
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.76'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.35%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("E36").Value = '  -4.40%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +28.80%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.14'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.48%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0473'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.27%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.75'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +36.83%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D41").Value = '0.0₃0701'
$ws.Range("E41").Value = '  -9.60%  '
$ws.Range("B42").Value = 'ThetaToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.91'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.10%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.142'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.80%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.65%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.22'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.29%  '
$ws.Range("B46").Value = 'TheGraph'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.326'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +8.49%  '
$ws.Range("B47").Value = 'WEMIXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.69'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +3.63%  '
$ws.Range("B48").Value = 'LidoDAOToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.32'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("E49").Value = '  -2.16%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '141.09'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.17%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.44%  '
